$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 : Réalisation / Réinitialisation de la carte SD ... ---
$ws.Range("A24").Value = "Réalisation"
$ws.Range("A24").WrapText = $true
$ws.Range("B24").Value = "Réinitialisation de la carte SD du raspberry Pi avec une nouvelle installation de Raspbian"
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 43509
$ws.Range("E24").Value = ""
$ws.Rows.Item(24).RowHeight = 45

# --- Row 25 : Conception / Mise à jour du MLD ... ---
$ws.Range("A25").Value = "Conception"
$ws.Range("A25").WrapText = $true
$ws.Range("B25").Value = "Mise à jour du MLD avec les informations récoltées depuis le cahier des charges, modification du MCD pour prendre en compte la table d'abonnement que peux obtenir un utilisateur"
$ws.Range("C25").Value = 1.5
$ws.Range("D25").Value = 43510
$ws.Range("E25").Value = "Je transfert le MLD que j'avais réalisé précedemment sur papier sur MySQL Workbench et j'y ajoute les champs auquel je n'avais pas pensé lors de la première conception du MLD et du MCD."
$ws.Rows.Item(25).RowHeight = 105

# --- Row 26 : Conception / Définition des activités possible ... ---
$ws.Range("A26").Value = "Conception"
$ws.Range("A26").WrapText = $true
$ws.Range("B26").Value = "Définition des activités possible, qui seront par la suite détaillée dans un diagramme UML"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 43510
$ws.Range("E26").Value = ""
$ws.Rows.Item(26).RowHeight = 60

# --- Row 27 : Analyse / Ajout d'endpoints à l'API ... ---
$ws.Range("A27").Value = "Analyse"
$ws.Range("A27").WrapText = $true
$ws.Range("B27").Value = "Ajout d'endpoints à l'API suite à la mise à jour du MLD"
$ws.Range("C27").Value = 0.5
$ws.Range("D27").Value = 43510
$ws.Range("E27").Value = ""
$ws.Rows.Item(27).RowHeight = 30

# --- View state: keep header frozen (1 row) and bring selection to E27 ---
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$aw.FreezePanes = $true
$ws.Range("E27").Select() | Out-Null

Write-Host "edit complete"
